# Updates the cryptos list (Price and Volume(1h) columns) to reflect the
# latest scrape, including a swap of the WrappedEther/ShibaInu rows (17/18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.012.94"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.587.95"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'527.03"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "'139.24"
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("D8").Value = "'0.564"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").Value = "2.596.87"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "'6.43"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  -3.35%  "
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "3.047.28"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "58.944.94"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "'20.53"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.574.24"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'344.37"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'10.07"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "'6.43"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'66.54"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'0.406"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'7.07"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "0.0₃0722"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").Value = "'5.90"
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("D33").Value = "'18.72"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'149.67"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'3.97"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").Value = "'36.80"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").Value = "'0.829"
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("D40").Value = "'0.811"
$ws.Range("E40").Value = "  -6.43%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "'270.02"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "'10.76"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").Value = "'18.38"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").Value = "1.963.13"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "'18.24"
$ws.Range("E51").Value = "  -2.58%  "
